$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.313.92"
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.391.12"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.77"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.79"
$ws.Range("E6").Value = "  +1.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("E9").Value = "  +8.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.590"
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.74"
$ws.Range("E11").Value = "  +3.68%  "
$ws.Range("E12").Value = "  +3.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "682.62"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.63"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.937.26"
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.388.17"
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.410.81"
$ws.Range("E17").Value = "  +1.94%  "
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("E20").Value = "  +2.81%  "
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.44"
$ws.Range("E22").Value = "  +1.33%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.78"
$ws.Range("E24").Value = "  +5.06%  "
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.73"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.25"
$ws.Range("E28").Value = "  +3.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.75"
$ws.Range("E29").Value = "  +1.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.01"
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("B31").Value = "dogwifhat"
$ws.Range("C31").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.71"
$ws.Range("E31").Value = "  +10.27%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.21"
$ws.Range("E32").Value = "  +1.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "555.04"
$ws.Range("E33").Value = "  -3.63%  "
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.14"
$ws.Range("E35").Value = "  +1.90%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.708.34"
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("E38").Value = "  +6.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.18"
$ws.Range("E39").Value = "  +1.33%  "
$ws.Range("E40").Value = "  +0.95%  "
$ws.Range("E41").Value = "  +3.87%  "
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("E44").Value = "  +4.14%  "
$ws.Range("E45").Value = "  -3.03%  "
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.40"
$ws.Range("E48").Value = "  +5.36%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.83"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("E51").Value = "  -1.90%  "
